$p = $ppt.ActivePresentation
Get-Member -InputObject $p.Slides
